$arr = New-Object "object[,]" 21,20
$arr[0,0] = 'Sending cluster'
$arr[0,1] = 'Ligand symbol'
$arr[0,2] = 'Receptor symbol'
$arr[0,3] = 'Target cluster'
$arr[0,4] = 'Ligand-expressing cells'
$arr[0,5] = 'Ligand detection rate'
$arr[0,6] = 'Ligand average expression value'
$arr[0,7] = 'Ligand total expression value'
$arr[0,8] = 'Ligand derived specificity of average expression value'
$arr[0,9] = 'Ligand derived specificity of total expression value'
$arr[0,10] = 'Receptor-expressing cells'
$arr[0,11] = 'Receptor detection rate'
$arr[0,12] = 'Receptor average expression value'
$arr[0,13] = 'Receptor total expression value'
$arr[0,14] = 'Receptor derived specificity of average expression value'
$arr[0,15] = 'Receptor derived specificity of total expression value'
$arr[0,16] = 'Edge average expression weight'
$arr[0,17] = 'Edge total expression weight'
$arr[0,18] = 'Edge average expression derived specificity'
$arr[0,19] = 'Edge total expression derived specificity'
$arr[1,0] = 'ECs'
$arr[1,1] = 'Col4a3'
$arr[1,2] = 'Cd93'
$arr[1,3] = 'ECs'
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 0.165709
$arr[1,7] = 0.497127
$arr[1,8] = 0.4546154542569759
$arr[1,9] = 0.4546154542569759
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 211.2725676666666
$arr[1,13] = 633.8177029999999
$arr[1,14] = 0.6324644927232657
$arr[1,15] = 0.6324644927232657
$arr[1,16] = 35.00976591547566
$arr[1,17] = 315.087893239281
$arr[1,18] = 0.2875281326607952
$arr[1,19] = 0.2875281326607952
$arr[2,0] = 'ECs'
$arr[2,1] = 'Col4a3'
$arr[2,2] = 'Cd93'
$arr[2,3] = 'Inflammatory-Mac'
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 0.165709
$arr[2,7] = 0.497127
$arr[2,8] = 0.4546154542569759
$arr[2,9] = 0.4546154542569759
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 59.36675400000001
$arr[2,13] = 178.100262
$arr[2,14] = 0.1777200152765546
$arr[2,15] = 0.1777200152765546
$arr[2,16] = 9.837605438586001
$arr[2,17] = 88.53844894727401
$arr[2,18] = 0.08079426547550757
$arr[2,19] = 0.08079426547550754
$arr[3,0] = 'ECs'
$arr[3,1] = 'Col4a3'
$arr[3,2] = 'Cd93'
$arr[3,3] = 'MuSCs'
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 0.165709
$arr[3,7] = 0.497127
$arr[3,8] = 0.4546154542569759
$arr[3,9] = 0.4546154542569759
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 0.4593846666666666
$arr[3,13] = 1.378154
$arr[3,14] = 0.001375211620595172
$arr[3,15] = 0.001375211620595172
$arr[3,16] = 0.07612417372866666
$arr[3,17] = 0.6851175635579999
$arr[3,18] = 0.0006251924555963462
$arr[3,19] = 0.0006251924555963461
$arr[4,0] = 'ECs'
$arr[4,1] = 'Col4a3'
$arr[4,2] = 'Cd93'
$arr[4,3] = 'Resolving-Mac'
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 0.165709
$arr[4,7] = 0.497127
$arr[4,8] = 0.4546154542569759
$arr[4,9] = 0.4546154542569759
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 62.94782133333333
$arr[4,13] = 188.843464
$arr[4,14] = 0.1884402803795846
$arr[4,15] = 0.1884402803795846
$arr[4,16] = 10.43102052532533
$arr[4,17] = 93.87918472792799
$arr[4,18] = 0.08566786366507677
$arr[4,19] = 0.08566786366507675
$arr[5,0] = 'FAPs'
$arr[5,1] = 'Col4a3'
$arr[5,2] = 'Cd93'
$arr[5,3] = 'ECs'
$arr[5,4] = 2
$arr[5,5] = 0.6666666666666666
$arr[5,6] = 0.1058106666666667
$arr[5,7] = 0.317432
$arr[5,8] = 0.2902869747080734
$arr[5,9] = 0.2902869747080733
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 211.2725676666666
$arr[5,13] = 633.8177029999999
$arr[5,14] = 0.6324644927232657
$arr[5,15] = 0.6324644927232657
$arr[5,16] = 22.35489123318844
$arr[5,17] = 201.194021098696
$arr[5,18] = 0.1835962042029131
$arr[5,19] = 0.183596204202913
$arr[6,0] = 'FAPs'
$arr[6,1] = 'Col4a3'
$arr[6,2] = 'Cd93'
$arr[6,3] = 'Inflammatory-Mac'
$arr[6,4] = 2
$arr[6,5] = 0.6666666666666666
$arr[6,6] = 0.1058106666666667
$arr[6,7] = 0.317432
$arr[6,8] = 0.2902869747080734
$arr[6,9] = 0.2902869747080733
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 59.36675400000001
$arr[6,13] = 178.100262
$arr[6,14] = 0.1777200152765546
$arr[6,15] = 0.1777200152765546
$arr[6,16] = 6.281635818576
$arr[6,17] = 56.534722367184
$arr[6,18] = 0.05158980557970361
$arr[6,19] = 0.05158980557970359
$arr[7,0] = 'FAPs'
$arr[7,1] = 'Col4a3'
$arr[7,2] = 'Cd93'
$arr[7,3] = 'MuSCs'
$arr[7,4] = 2
$arr[7,5] = 0.6666666666666666
$arr[7,6] = 0.1058106666666667
$arr[7,7] = 0.317432
$arr[7,8] = 0.2902869747080734
$arr[7,9] = 0.2902869747080733
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 0.4593846666666666
$arr[7,13] = 1.378154
$arr[7,14] = 0.001375211620595172
$arr[7,15] = 0.001375211620595172
$arr[7,16] = 0.04860779783644444
$arr[7,17] = 0.437470180528
$arr[7,18] = 0.0003992060209259593
$arr[7,19] = 0.0003992060209259592
$arr[8,0] = 'FAPs'
$arr[8,1] = 'Col4a3'
$arr[8,2] = 'Cd93'
$arr[8,3] = 'Resolving-Mac'
$arr[8,4] = 2
$arr[8,5] = 0.6666666666666666
$arr[8,6] = 0.1058106666666667
$arr[8,7] = 0.317432
$arr[8,8] = 0.2902869747080734
$arr[8,9] = 0.2902869747080733
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 62.94782133333333
$arr[8,13] = 188.843464
$arr[8,14] = 0.1884402803795846
$arr[8,15] = 0.1884402803795846
$arr[8,16] = 6.660550940494222
$arr[8,17] = 59.94495846444799
$arr[8,18] = 0.05470175890453073
$arr[8,19] = 0.05470175890453072
$arr[9,0] = 'Inflammatory-Mac'
$arr[9,1] = 'Col4a3'
$arr[9,2] = 'Cd93'
$arr[9,3] = 'ECs'
$arr[9,4] = 1
$arr[9,5] = 0.3333333333333333
$arr[9,6] = 0.001077
$arr[9,7] = 0.003231
$arr[9,8] = 0.002954702787626279
$arr[9,9] = 0.002954702787626278
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 211.2725676666666
$arr[9,13] = 633.8177029999999
$arr[9,14] = 0.6324644927232657
$arr[9,15] = 0.6324644927232657
$arr[9,16] = 0.227540555377
$arr[9,17] = 2.047864998393
$arr[9,18] = 0.001868744599724074
$arr[9,19] = 0.001868744599724073
$arr[10,0] = 'Inflammatory-Mac'
$arr[10,1] = 'Col4a3'
$arr[10,2] = 'Cd93'
$arr[10,3] = 'Inflammatory-Mac'
$arr[10,4] = 1
$arr[10,5] = 0.3333333333333333
$arr[10,6] = 0.001077
$arr[10,7] = 0.003231
$arr[10,8] = 0.002954702787626279
$arr[10,9] = 0.002954702787626278
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 59.36675400000001
$arr[10,13] = 178.100262
$arr[10,14] = 0.1777200152765546
$arr[10,15] = 0.1777200152765546
$arr[10,16] = 0.06393799405800001
$arr[10,17] = 0.575441946522
$arr[10,18] = 0.0005251098245546208
$arr[10,19] = 0.0005251098245546206
$arr[11,0] = 'Inflammatory-Mac'
$arr[11,1] = 'Col4a3'
$arr[11,2] = 'Cd93'
$arr[11,3] = 'MuSCs'
$arr[11,4] = 1
$arr[11,5] = 0.3333333333333333
$arr[11,6] = 0.001077
$arr[11,7] = 0.003231
$arr[11,8] = 0.002954702787626279
$arr[11,9] = 0.002954702787626278
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 0.4593846666666666
$arr[11,13] = 1.378154
$arr[11,14] = 0.001375211620595172
$arr[11,15] = 0.001375211620595172
$arr[11,16] = 0.0004947572859999999
$arr[11,17] = 0.004452815573999999
$arr[11,18] = [double]"4.063341608948608E-06"
$arr[11,19] = [double]"4.063341608948607E-06"
$arr[12,0] = 'Inflammatory-Mac'
$arr[12,1] = 'Col4a3'
$arr[12,2] = 'Cd93'
$arr[12,3] = 'Resolving-Mac'
$arr[12,4] = 1
$arr[12,5] = 0.3333333333333333
$arr[12,6] = 0.001077
$arr[12,7] = 0.003231
$arr[12,8] = 0.002954702787626279
$arr[12,9] = 0.002954702787626278
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 62.94782133333333
$arr[12,13] = 188.843464
$arr[12,14] = 0.1884402803795846
$arr[12,15] = 0.1884402803795846
$arr[12,16] = 0.067794803576
$arr[12,17] = 0.610153232184
$arr[12,18] = 0.0005567850217386364
$arr[12,19] = 0.0005567850217386361
$arr[13,0] = 'MuSCs'
$arr[13,1] = 'Col4a3'
$arr[13,2] = 'Cd93'
$arr[13,3] = 'ECs'
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 0.07261866666666666
$arr[13,7] = 0.217856
$arr[13,8] = 0.1992261623339866
$arr[13,9] = 0.1992261623339865
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 211.2725676666666
$arr[13,13] = 633.8177029999999
$arr[13,14] = 0.6324644927232657
$arr[13,15] = 0.6324644927232657
$arr[13,16] = 15.34233216719644
$arr[13,17] = 138.080989504768
$arr[13,18] = 0.1260034736977678
$arr[13,19] = 0.1260034736977678
$arr[14,0] = 'MuSCs'
$arr[14,1] = 'Col4a3'
$arr[14,2] = 'Cd93'
$arr[14,3] = 'Inflammatory-Mac'
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 0.07261866666666666
$arr[14,7] = 0.217856
$arr[14,8] = 0.1992261623339866
$arr[14,9] = 0.1992261623339865
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 59.36675400000001
$arr[14,13] = 178.100262
$arr[14,14] = 0.1777200152765546
$arr[14,15] = 0.1777200152765546
$arr[14,16] = 4.311134519808
$arr[14,17] = 38.800210678272
$arr[14,18] = 0.03540647661348544
$arr[14,19] = 0.03540647661348543
$arr[15,0] = 'MuSCs'
$arr[15,1] = 'Col4a3'
$arr[15,2] = 'Cd93'
$arr[15,3] = 'MuSCs'
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 0.07261866666666666
$arr[15,7] = 0.217856
$arr[15,8] = 0.1992261623339866
$arr[15,9] = 0.1992261623339865
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 0.4593846666666666
$arr[15,13] = 1.378154
$arr[15,14] = 0.001375211620595172
$arr[15,15] = 0.001375211620595172
$arr[15,16] = 0.03335990198044444
$arr[15,17] = 0.300239117824
$arr[15,18] = 0.0002739781335682785
$arr[15,19] = 0.0002739781335682785
$arr[16,0] = 'MuSCs'
$arr[16,1] = 'Col4a3'
$arr[16,2] = 'Cd93'
$arr[16,3] = 'Resolving-Mac'
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 0.07261866666666666
$arr[16,7] = 0.217856
$arr[16,8] = 0.1992261623339866
$arr[16,9] = 0.1992261623339865
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = 62.94782133333333
$arr[16,13] = 188.843464
$arr[16,14] = 0.1884402803795846
$arr[16,15] = 0.1884402803795846
$arr[16,16] = 4.571186854798222
$arr[16,17] = 41.14068169318399
$arr[16,18] = 0.03754223388916507
$arr[16,19] = 0.03754223388916506
$arr[17,0] = 'Resolving-Mac'
$arr[17,1] = 'Col4a3'
$arr[17,2] = 'Cd93'
$arr[17,3] = 'ECs'
$arr[17,4] = 1
$arr[17,5] = 0.3333333333333333
$arr[17,6] = 0.01928833333333333
$arr[17,7] = 0.057865
$arr[17,8] = 0.05291670591333787
$arr[17,9] = 0.05291670591333786
$arr[17,10] = 3
$arr[17,11] = 1
$arr[17,12] = 211.2725676666666
$arr[17,13] = 633.8177029999999
$arr[17,14] = 0.6324644927232657
$arr[17,15] = 0.6324644927232657
$arr[17,16] = 4.075095709343889
$arr[17,17] = 36.67586138409499
$arr[17,18] = 0.03346793756206547
$arr[17,19] = 0.03346793756206546
$arr[18,0] = 'Resolving-Mac'
$arr[18,1] = 'Col4a3'
$arr[18,2] = 'Cd93'
$arr[18,3] = 'Inflammatory-Mac'
$arr[18,4] = 1
$arr[18,5] = 0.3333333333333333
$arr[18,6] = 0.01928833333333333
$arr[18,7] = 0.057865
$arr[18,8] = 0.05291670591333787
$arr[18,9] = 0.05291670591333786
$arr[18,10] = 3
$arr[18,11] = 1
$arr[18,12] = 59.36675400000001
$arr[18,13] = 178.100262
$arr[18,14] = 0.1777200152765546
$arr[18,15] = 0.1777200152765546
$arr[18,16] = 1.14508574007
$arr[18,17] = 10.30577166063
$arr[18,18] = 0.009404357783303351
$arr[18,19] = 0.00940435778330335
$arr[19,0] = 'Resolving-Mac'
$arr[19,1] = 'Col4a3'
$arr[19,2] = 'Cd93'
$arr[19,3] = 'MuSCs'
$arr[19,4] = 1
$arr[19,5] = 0.3333333333333333
$arr[19,6] = 0.01928833333333333
$arr[19,7] = 0.057865
$arr[19,8] = 0.05291670591333787
$arr[19,9] = 0.05291670591333786
$arr[19,10] = 3
$arr[19,11] = 1
$arr[19,12] = 0.4593846666666666
$arr[19,13] = 1.378154
$arr[19,14] = 0.001375211620595172
$arr[19,15] = 0.001375211620595172
$arr[19,16] = 0.008860764578888888
$arr[19,17] = 0.07974688121
$arr[19,18] = [double]"7.277166889563949E-05"
$arr[19,19] = [double]"7.277166889563948E-05"
$arr[20,0] = 'Resolving-Mac'
$arr[20,1] = 'Col4a3'
$arr[20,2] = 'Cd93'
$arr[20,3] = 'Resolving-Mac'
$arr[20,4] = 1
$arr[20,5] = 0.3333333333333333
$arr[20,6] = 0.01928833333333333
$arr[20,7] = 0.057865
$arr[20,8] = 0.05291670591333787
$arr[20,9] = 0.05291670591333786
$arr[20,10] = 3
$arr[20,11] = 1
$arr[20,12] = 62.94782133333333
$arr[20,13] = 188.843464
$arr[20,14] = 0.1884402803795846
$arr[20,15] = 0.1884402803795846
$arr[20,16] = 1.214158560484444
$arr[20,17] = 10.92742704436
$arr[20,18] = 0.009971638899073411
$arr[20,19] = 0.009971638899073409

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1:T21").Value = $arr